# ReportingOrganisationGroup.xlsx: columns D ("codeforiati:group-name")
# and E ("codeforiati:group-code") were swapped for every row (including
# the header row) so that D now holds the group CODE and E holds the
# group NAME.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
